# Update column E: header and values to represent "Average Height (mm)"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header in E1
$ws.Range("E1").Value = "Average Height (mm)"

# Compute the average of the Height (mm) column (B2:B101)
$avgHeight = $ws.Application.WorksheetFunction.Average($ws.Range("B2:B101"))

# Fill E2:E101 with the computed average value (constant values, not formulas)
$ws.Range("E2:E101").Value = $avgHeight
